$wb = $excel.ActiveWorkbook

# --- Queries sheet: update the SQL text (adds ORDER BY clause) ---
$wsQueries = $wb.Worksheets.Item("Queries")
$queryText = "Select IMU_MENU_ID as 'Menu ID'`n  ,IVN_VDN_NUM as VDN`n  ,IVN_VDN_DESC as 'VDN Description'`n  ,IVN_VIP_VDN_NUM as 'VIP VDN'`n  ,IVN_MENU_OPT as 'Option'`n  ,IVN_SESS_DNIS as DNIS`n  FROM [IVR_VDN_NUM] Order By IVN_VDN_DESC Asc;"
$wsQueries.Range("A2").Value = $queryText

# --- Create sheet: replace test data values (leading ' preserves the
#     existing quotePrefix text style instead of minting a new style) ---
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("C3").Value = "'4356"
$wsCreate.Range("E3").Value = "'4728"
$wsCreate.Range("F3").Value = "'Rest"
$wsCreate.Range("F4").Value = "'User"
$wsCreate.Range("F5").Value = "'Data"
$wsCreate.Range("E4").Value = "'5267"

# --- Update selections on each sheet (restore Edit sheet as active tab last) ---
$wsCreate.Range("B7").Select()

$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("F6").Select()

$wsQueries.Range("A3").Select()

$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("H3").Select()
